# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting from the neighboring header cell (G1)
# so it picks up the same bold/border/alignment style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Per-row "Save" flag values for H2:H29.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
}

foreach ($row in 2..29) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
